$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force literal text (avoid Excel's numeric auto-coercion for
    # number-like strings), then strip the resulting quote-prefix
    # style so the cell's style index is left untouched.
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.352.08"
Set-TextValue $ws.Range("E2") "  -0.69%  "
Set-TextValue $ws.Range("D3") "1.713.89"
Set-TextValue $ws.Range("E3") "  -1.35%  "
Set-TextValue $ws.Range("D4") "0.9961"
Set-TextValue $ws.Range("E4") "  -0.35%  "
Set-TextValue $ws.Range("D5") "240.48"
Set-TextValue $ws.Range("E5") "  -2.57%  "
Set-TextValue $ws.Range("D6") "0.9970"
Set-TextValue $ws.Range("E6") "  -0.31%  "
Set-TextValue $ws.Range("D7") "0.4861"
Set-TextValue $ws.Range("E7") "  -1.30%  "
Set-TextValue $ws.Range("D8") "0.2583"
Set-TextValue $ws.Range("E8") "  -3.17%  "
Set-TextValue $ws.Range("D9") "0.06174"
Set-TextValue $ws.Range("E9") "  -1.94%  "
Set-TextValue $ws.Range("D10") "1.713.47"
Set-TextValue $ws.Range("E10") "  -1.33%  "
Set-TextValue $ws.Range("D11") "0.06947"
Set-TextValue $ws.Range("E11") "  -1.51%  "
Set-TextValue $ws.Range("D12") "15.47"
Set-TextValue $ws.Range("E12") "  -1.57%  "
Set-TextValue $ws.Range("D13") "4.467"
Set-TextValue $ws.Range("E13") "  -2.96%  "
Set-TextValue $ws.Range("D14") "0.5968"
Set-TextValue $ws.Range("E14") "  -2.46%  "
Set-TextValue $ws.Range("D15") "76.41"
Set-TextValue $ws.Range("E15") "  -1.41%  "
Set-TextValue $ws.Range("D16") "0.9972"
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D17") "26.242.88"
Set-TextValue $ws.Range("E17") "  -1.09%  "
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D18") "0.9963"
Set-TextValue $ws.Range("E18") "  -0.39%  "
Set-TextValue $ws.Range("D19") "0.000007105"
Set-TextValue $ws.Range("E19") "  -3.63%  "
Set-TextValue $ws.Range("D20") "11.22"
Set-TextValue $ws.Range("E20") "  -2.81%  "
Set-TextValue $ws.Range("D21") "1.932.97"
Set-TextValue $ws.Range("E21") "  -1.30%  "
Set-TextValue $ws.Range("D22") "4.404"
Set-TextValue $ws.Range("E22") "  -4.12%  "
Set-TextValue $ws.Range("D23") "8.423"
Set-TextValue $ws.Range("E23") "  -3.27%  "
Set-TextValue $ws.Range("D24") "5.041"
Set-TextValue $ws.Range("E24") "  -4.14%  "
Set-TextValue $ws.Range("D25") "136.07"
Set-TextValue $ws.Range("E25") "  -2.68%  "
Set-TextValue $ws.Range("D26") "15.15"
Set-TextValue $ws.Range("E26") "  -2.02%  "
Set-TextValue $ws.Range("D27") "1.395"
Set-TextValue $ws.Range("E27") "  -1.78%  "
Set-TextValue $ws.Range("D28") "1.730"
Set-TextValue $ws.Range("E28") "  -1.94%  "
Set-TextValue $ws.Range("D29") "105.28"
Set-TextValue $ws.Range("E29") "  -2.60%  "
Set-TextValue $ws.Range("D30") "3.872"
Set-TextValue $ws.Range("E30") "  -4.16%  "
Set-TextValue $ws.Range("E31") "  -1.55%  "
Set-TextValue $ws.Range("D32") "3.604"
Set-TextValue $ws.Range("D33") "0.04421"
Set-TextValue $ws.Range("D34") "2.597"
Set-TextValue $ws.Range("E34") "  -0.53%  "
Set-TextValue $ws.Range("D35") "0.9896"
Set-TextValue $ws.Range("E35") "  -2.04%  "
Set-TextValue $ws.Range("D36") "0.6167"
Set-TextValue $ws.Range("E36") "  -3.20%  "
Set-TextValue $ws.Range("D37") "0.9351"
Set-TextValue $ws.Range("E37") "  +4.35%  "
Set-TextValue $ws.Range("D38") "1.975"
Set-TextValue $ws.Range("E38") "  -1.94%  "
Set-TextValue $ws.Range("D39") "2.369"
Set-TextValue $ws.Range("E39") "  -1.50%  "
Set-TextValue $ws.Range("D40") "0.9962"
Set-TextValue $ws.Range("E40") "  -0.66%  "
Set-TextValue $ws.Range("D41") "0.01472"
Set-TextValue $ws.Range("E41") "  -2.66%  "
Set-TextValue $ws.Range("D42") "99.61"
Set-TextValue $ws.Range("E42") "  -2.55%  "
Set-TextValue $ws.Range("D43") "5.376"
Set-TextValue $ws.Range("E43") "  -0.50%  "
Set-TextValue $ws.Range("D44") "0.3799"
Set-TextValue $ws.Range("E44") "  -2.82%  "
Set-TextValue $ws.Range("D45") "6.825"
Set-TextValue $ws.Range("E45") "  -0.86%  "
Set-TextValue $ws.Range("D46") "0.1148"
Set-TextValue $ws.Range("E46") "  -3.34%  "
Set-TextValue $ws.Range("E47") "  -1.07%  "
Set-TextValue $ws.Range("D48") "30.59"
Set-TextValue $ws.Range("E48") "  +0.11%  "
Set-TextValue $ws.Range("D49") "7.703"
Set-TextValue $ws.Range("E49") "  -1.10%  "
Set-TextValue $ws.Range("D50") "51.02"
Set-TextValue $ws.Range("E50") "  -1.60%  "
Set-TextValue $ws.Range("D51") "1.210"
Set-TextValue $ws.Range("E51") "  -4.69%  "
